# Apply the "week-1.en.md_word.pptx" slide2 download-links edit:
#   PPTX-MD (rId4)              -> DOCX-MD (keeps rId4)
#   (new run)                   -> PPTX-MD (reuses rId5's existing target)
#   PPTX-MS (was rId5)          -> PPTX-MS (now rId6, a brand new relationship)
#
# i.e. "Download PDF-MS, PDF-MD, PPTX-MD, PPTX-MS"
#  ->  "Download PDF-MS, PDF-MD, DOCX-MD, PPTX-MD, PPTX-MS"

$p = $ppt.ActivePresentation

# Locate the slide/shape that holds the "Download ..." links line instead of
# hard-coding indices, so the script is resilient to reordering.
$targetSlide = $null
$targetShape = $null
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTextFrame) {
            $t = $shape.TextFrame.TextRange.Text
            if ($t.IndexOf("PPTX-MD") -ge 0) {
                $targetSlide = $slide
                $targetShape = $shape
            }
        }
    }
}

$tr = $targetShape.TextFrame.TextRange

# Step 1: rename the existing "PPTX-MD" run to "DOCX-MD" in place. Replacing
# the text of an exact-run character range keeps the run's own formatting
# (its <a:hlinkClick r:id="rId4"/>) untouched - only the <a:t> changes.
$full = $tr.Text
$idx = $full.IndexOf("PPTX-MD")
$run = $tr.Characters($idx + 1, 7)
$run.Text = "DOCX-MD"

# Step 2: grow the plain ", " separator that precedes "PPTX-MS" into
# ", PPTX-MD, " (still a single, unformatted run at this point).
$full = $tr.Text
$msIdx = $full.IndexOf("PPTX-MS")
$sepStart = $msIdx - 2 + 1
$sep = $tr.Characters($sepStart, 2)
$sep.Text = ", PPTX-MD, "

# Step 3: turn the newly inserted "PPTX-MD" word into a hyperlink that
# targets the same file the original "PPTX-MS" run (rId5) pointed to - the
# host package reuses that existing relationship instead of cloning it.
$full = $tr.Text
$mdIdx = $full.IndexOf("PPTX-MD")
$mdRun = $tr.Characters($mdIdx + 1, 7)
$mdRun.ActionSettings(1).Hyperlink.Address = "1_veri_tabani_giris_week_1.pptx"

# Step 4: give the trailing "PPTX-MS" its own, brand new hyperlink
# relationship (rId6), leaving rId4/rId5 exactly as they were.
$full = $tr.Text
$msIdx2 = $full.LastIndexOf("PPTX-MS")
$msRun = $tr.Characters($msIdx2 + 1, 7)
$msRun.ActionSettings(1).Hyperlink.Address = "week-1.en.md_word.pptx"

Write-Host ("Updated text: " + $tr.Text)
